$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-31: column A is present but blank, column B = "yes".
# Setting .Style after the blank .Value keeps the (empty) cell materialised
# in the sheet (instead of being dropped as a totally absent cell) without
# pulling in a new number-format / style definition.
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = ""
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = "yes"
}

# Rows 32-36: column A holds specific Brnum values, column B = "yes".
$brnums = @("BR50041", "BR50042", "BR50045", "BR50047", "BR50050")
$row = 32
foreach ($brnum in $brnums) {
    $ws.Cells.Item($row, 1).Value = $brnum
    $ws.Cells.Item($row, 2).Value = "yes"
    $row++
}
